$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Level 10 unlock gains "Absorb Hex" (merged in from the old C19 cell)
$ws.Range("B12").Value = "Light Poison Resistant Aliens, Absorb Hex"

# "Upcoming Features" column (C) entries for Absorb Hex / Cursed Rocket / Plasma Hex
# are cleared out - those features shipped, so the column no longer lists them here.
$ws.Range("C19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("C25").Value = ""

# "Omega Rocket" moves up from B26 into B25 (row 24 now has no Unlock(s) entry)
# Grab the existing row-25 cell format (bordered data-row style) before writing
# into the previously-empty B25 cell, so it matches its neighbours.
$ws.Range("A25").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = "Omega Rocket"
$ws.Range("B26").Value = ""

# Restore the cursor/selection to match the saved workbook state
$ws.Range("B13").Select()
